$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended to the EUR->ARS rate log.
$row = 11

# Force text storage (not Excel's date/time auto-conversion) while keeping
# the cells on the sheet's default "Normal" style, matching the rest of the
# table (no explicit per-cell style).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-08"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "15:20:28"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1616.0155 ARS"
$ws.Cells.Item($row, 3).Style = "Normal"
